$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for column A (SITC1 -> SITC_1)
$ws.Range("A1").Value = "SITC_1"

# Strip the leading "N." numeric prefix from each classification description
$ws.Range("B2").Value = "Food & Live Animals"
$ws.Range("B3").Value = "Beverages & Tobacco"
$ws.Range("B4").Value = "Crude Materials Except Fuels"
$ws.Range("B5").Value = "Mineral Fuels"
$ws.Range("B6").Value = "Oils, Fats & Waxes"
$ws.Range("B7").Value = "Chemical Products"
$ws.Range("B8").Value = "Basic Manufactured Products"
$ws.Range("B9").Value = "Machines & Transport Equipment"
$ws.Range("B10").Value = "Miscellaneous Manufactured Goods"
$ws.Range("B11").Value = "Goods N.E.C."

# Move the active selection to match the author's final cursor position
$ws.Range("B15").Select()
